$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' (before) -> 'Bitcoin' (after)
$ws.Range("D2").Value = "26.597.17"
$ws.Range("E2").Value = "  +1.23%  "

# Row 3: 'Ethereum' (before) -> 'Ethereum' (after)
$ws.Range("D3").Value = "1.633.80"
$ws.Range("E3").Value = "  +0.54%  "

# Row 4: 'TetherUSD' (before) -> 'TetherUSD' (after)
$ws.Range("E4").Value = "  -0.16%  "

# Row 5: 'BNB' (before) -> 'BNB' (after)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.11"
$ws.Range("E5").Value = "  +0.12%  "

# Row 6: 'XRP' (before) -> 'XRP' (after)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("E6").Value = "  +2.58%  "

# Row 7: 'USDC' (before) -> 'USDC' (after)
$ws.Range("E7").Value = "  -0.13%  "

# Row 8: 'Cardano' (before) -> 'Cardano' (after)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +2.18%  "

# Row 9: 'Dogecoin' (before) -> 'Dogecoin' (after)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0625"
$ws.Range("E9").Value = "  +1.52%  "

# Row 10: 'Solana' (before) -> 'Solana' (after)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.17"
$ws.Range("E10").Value = "  +1.30%  "

# Row 11: 'TRON' (before) -> 'TRON' (after)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0840"
$ws.Range("E11").Value = "  +2.62%  "

# Row 12: 'WrappedliquidstakedEther2.0' (before) -> 'WrappedliquidstakedEther2.0' (after)
$ws.Range("D12").Value = "1.857.54"
$ws.Range("E12").Value = "  +0.25%  "

# Row 13: 'WrappedEther' (before) -> 'WrappedEther' (after)
$ws.Range("D13").Value = "1.633.61"
$ws.Range("E13").Value = "  +0.48%  "

# Row 14: 'Polkadot' (before) -> 'Polkadot' (after)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("E14").Value = "  +1.64%  "

# Row 15: 'Polygon' (before) -> 'Polygon' (after)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  +1.63%  "

# Row 16: 'WrappedBTC' (before) -> 'WrappedBTC' (after)
$ws.Range("D16").Value = "26.580.00"
$ws.Range("E16").Value = "  +1.11%  "

# Row 17: 'Litecoin' (before) -> 'Litecoin' (after)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.35"
$ws.Range("E17").Value = "  +1.34%  "

# Row 18: 'ShibaInu' (before) -> 'ShibaInu' (after)
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("E18").Value = "  +2.04%  "

# Row 19: 'BitcoinCash' (before) -> 'BitcoinCash' (after)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.29"
$ws.Range("E19").Value = "  +8.09%  "

# Row 20: 'Dai' (before) -> 'Dai' (after)
$ws.Range("E20").Value = "  -0.07%  "

# Row 21: 'Uniswap' (before) -> 'Uniswap' (after)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.31"
$ws.Range("E21").Value = "  +0.21%  "

# Row 22: 'Avalanche' (before) -> 'Avalanche' (after)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.51"
$ws.Range("E22").Value = "  +1.46%  "

# Row 23: 'Chainlink' (before) -> 'Chainlink' (after)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.21"
$ws.Range("E23").Value = "  +2.76%  "

# Row 24: 'Toncoin' (before) -> 'Toncoin' (after)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  -0.65%  "

# Row 25: 'Monero' (before) -> 'Monero' (after)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.19"
$ws.Range("E25").Value = "  +4.07%  "

# Row 26: 'BinanceUSD' (before) -> 'BinanceUSD' (after)
$ws.Range("E26").Value = "  -0.08%  "

# Row 27: 'Stellar' (before) -> 'Stellar' (after)
$ws.Range("E27").Value = "  +1.28%  "

# Row 28: 'Cosmos' (before) -> 'Cosmos' (after)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.88"
$ws.Range("E28").Value = "  +4.74%  "

# Row 29: 'EthereumClassic' (before) -> 'EthereumClassic' (after)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.53"
$ws.Range("E29").Value = "  +1.41%  "

# Row 30: 'Hedera' (before) -> 'Hedera' (after)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("E30").Value = "  -1.47%  "

# Row 31: 'PancakeSwap' (before) -> 'PancakeSwap' (after)
$ws.Range("E31").Value = "  -0.74%  "

# Row 32: 'Filecoin' (before) -> 'Filecoin' (after)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("E32").Value = "  +2.62%  "

# Row 33: 'InternetComputer(DFINITY)' (before) -> 'InternetComputer(DFINITY)' (after)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("E33").Value = "  -1.28%  "

# Row 34: 'LidoDAOToken' (before) -> 'LidoDAOToken' (after)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").Value = "  -0.50%  "

# Row 35: 'HuobiToken' (before) -> 'HuobiToken' (after)
$ws.Range("E35").Value = "  -0.62%  "

# Row 36: 'VeChain' (before) -> 'VeChain' (after)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0173"
$ws.Range("E36").Value = "  +4.33%  "

# Row 37: 'Maker' (before) -> 'Maker' (after)
$ws.Range("D37").Value = "1.177.79"
$ws.Range("E37").Value = "  +0.68%  "

# Row 38: 'ARBITRUM' (before) -> 'ARBITRUM' (after)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.807"
$ws.Range("E38").Value = "  +0.05%  "

# Row 40: 'ImmutableX' (before) -> 'ImmutableX' (after)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.507"
$ws.Range("E40").Value = "  +1.77%  "

# Row 41: 'MXToken' (before) -> 'MXToken' (after)
$ws.Range("E41").Value = "  -0.44%  "

# Row 42: 'TrustWalletToken' (before) -> 'FraxShare' (after)
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("E42").Value = "  +1.63%  "

# Row 43: 'FraxShare' (before) -> 'TrustWalletToken' (after)
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.791"
$ws.Range("E43").Value = "  -0.30%  "

# Row 44: 'RocketPoolETH' (before) -> 'RocketPoolETH' (after)
$ws.Range("D44").Value = "1.767.13"
$ws.Range("E44").Value = "  +0.24%  "

# Row 45: 'Quant' (before) -> 'Quant' (after)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.81"
$ws.Range("E45").Value = "  -0.64%  "

# Row 46: 'RenderToken' (before) -> 'RenderToken' (after)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.54"
$ws.Range("E46").Value = "  +1.29%  "

# Row 47: 'Aave' (before) -> 'Aave' (after)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.80"
$ws.Range("E47").Value = "  +1.04%  "

# Row 48: 'Cronos' (before) -> 'Cronos' (after)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  +0.80%  "

# Row 49: 'EnergySwap' (before) -> 'EnergySwap' (after)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.64"
$ws.Range("E49").Value = "  +4.68%  "

# Row 50: 'Mantle' (before) -> 'Mantle' (after)
$ws.Range("E50").Value = "  +0.09%  "

# Row 51: 'USDD' (before) -> 'USDD' (after)
$ws.Range("E51").Value = "  -0.19%  "
